$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Column A holds a date-like label ("01-07-2021") that must stay as plain text,
# matching the other cells in column A which are shared-string text values.
# Force text formatting before assignment so Excel does not auto-convert it
# to a date serial number, then restore the default (unstyled) appearance so
# the new row matches the look of the existing data rows.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "01-07-2021"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 4.5
$ws.Cells.Item($row, 3).Value = 3.6
$ws.Cells.Item($row, 4).Value = 5.9
$ws.Cells.Item($row, 5).Value = 3
$ws.Cells.Item($row, 6).Value = 6.7
